# Scheduled market-data refresh for the Zodiark_Profits workbook.
# Updates the FFXIV Market Board snapshot columns (H:N) -- currentAveragePrice,
# currentAveragePriceNQ/HQ, LevePriceNQ/HQ and LeveProfitNQ/HQ -- for the leve
# rows whose item prices moved since the last run, one job sheet at a time.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 2789
$ws.Range("I18").Value = 3250.8
$ws.Range("K18").Value = 3250.8
$ws.Range("M18").Value = -2966.8
$ws.Range("H19").Value = 3497.375
$ws.Range("I19").Value = 1797.8
$ws.Range("J19").Value = 6330
$ws.Range("K19").Value = 1797.8
$ws.Range("L19").Value = 6330
$ws.Range("M19").Value = -1622.8
$ws.Range("N19").Value = -6680
$ws.Range("H28").Value = 589117.7
$ws.Range("I28").Value = 1000376.1
$ws.Range("K28").Value = 1000376.1
$ws.Range("M28").Value = -999891.1
$ws.Range("H64").Value = 7044.75
$ws.Range("I64").Value = 6848
$ws.Range("J64").Value = 7438.25
$ws.Range("K64").Value = 6848
$ws.Range("L64").Value = 7438.25
$ws.Range("M64").Value = -6600
$ws.Range("N64").Value = -7934.25
$ws.Range("H67").Value = 7044.75
$ws.Range("I67").Value = 6848
$ws.Range("J67").Value = 7438.25
$ws.Range("K67").Value = 6848
$ws.Range("L67").Value = 7438.25
$ws.Range("M67").Value = -5990
$ws.Range("N67").Value = -9154.25
$ws.Range("H74").Value = 6924.5
$ws.Range("J74").Value = 9000
$ws.Range("L74").Value = 9000
$ws.Range("N74").Value = -10872
$ws.Range("H77").Value = 6924.5
$ws.Range("J77").Value = 9000
$ws.Range("L77").Value = 45000
$ws.Range("N77").Value = -54360
$ws.Range("H108").Value = 60000
$ws.Range("J108").Value = 60000
$ws.Range("L108").Value = 60000
$ws.Range("N108").Value = -67680
$ws.Range("H114").Value = 66000
$ws.Range("J114").Value = 66000
$ws.Range("L114").Value = 66000
$ws.Range("N114").Value = -74678
$ws.Range("H132").Value = 3129.3
$ws.Range("I132").Value = 2921.5557
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 8764.667099999999
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -6234.667099999999
$ws.Range("N132").Value = -20057

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2240
$ws.Range("I74").Value = 1552.258
$ws.Range("K74").Value = 1552.258
$ws.Range("M74").Value = -678.258
$ws.Range("H77").Value = 2240
$ws.Range("I77").Value = 1552.258
$ws.Range("K77").Value = 7761.29
$ws.Range("M77").Value = -3393.29
$ws.Range("H97").Value = 723.82355
$ws.Range("I97").Value = 496.25
$ws.Range("J97").Value = 926.1111
$ws.Range("K97").Value = 496.25
$ws.Range("L97").Value = 926.1111
$ws.Range("M97").Value = -0.25
$ws.Range("N97").Value = -1918.1111
$ws.Range("H122").Value = 6395.905
$ws.Range("I122").Value = 7402.5
$ws.Range("J122").Value = 4382.7144
$ws.Range("K122").Value = 22207.5
$ws.Range("L122").Value = 13148.1432
$ws.Range("M122").Value = -19757.5
$ws.Range("N122").Value = -18048.1432
$ws.Range("H132").Value = 5035.08
$ws.Range("I132").Value = 5551.3335
$ws.Range("J132").Value = 2324.75
$ws.Range("K132").Value = 16654.0005
$ws.Range("L132").Value = 6974.25
$ws.Range("M132").Value = -14124.0005
$ws.Range("N132").Value = -12034.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H57").Value = 100000
$ws.Range("J57").Value = 100000
$ws.Range("L57").Value = 100000
$ws.Range("N57").Value = -101440
$ws.Range("H122").Value = 50000
$ws.Range("J122").Value = 50000
$ws.Range("L122").Value = 50000
$ws.Range("N122").Value = -59800
$ws.Range("H136").Value = 100000
$ws.Range("J136").Value = 100000
$ws.Range("L136").Value = 100000
$ws.Range("N136").Value = -110200

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2424.375
$ws.Range("J31").Value = 2784.111
$ws.Range("L31").Value = 2784.111
$ws.Range("N31").Value = -3374.111
$ws.Range("H34").Value = 2424.375
$ws.Range("J34").Value = 2784.111
$ws.Range("L34").Value = 2784.111
$ws.Range("N34").Value = -3188.111
$ws.Range("H58").Value = 6614.231
$ws.Range("I58").Value = 3075.8
$ws.Range("K58").Value = 3075.8
$ws.Range("M58").Value = -2872.8
$ws.Range("H62").Value = 7163.125
$ws.Range("I62").Value = 4551.6665
$ws.Range("J62").Value = 14997.5
$ws.Range("K62").Value = 4551.6665
$ws.Range("L62").Value = 14997.5
$ws.Range("M62").Value = -3927.6665
$ws.Range("N62").Value = -16245.5
$ws.Range("H65").Value = 7163.125
$ws.Range("I65").Value = 4551.6665
$ws.Range("J65").Value = 14997.5
$ws.Range("K65").Value = 22758.3325
$ws.Range("L65").Value = 74987.5
$ws.Range("M65").Value = -19638.3325
$ws.Range("N65").Value = -81227.5
$ws.Range("H99").Value = 2099.0312
$ws.Range("I99").Value = 2021.1538
$ws.Range("K99").Value = 2021.1538
$ws.Range("M99").Value = -523.1538
$ws.Range("H102").Value = 57500
$ws.Range("J102").Value = 57500
$ws.Range("L102").Value = 57500
$ws.Range("N102").Value = -62368
$ws.Range("H109").Value = 31814.5
$ws.Range("J109").Value = 28999.666
$ws.Range("L109").Value = 28999.666
$ws.Range("N109").Value = -31079.666
$ws.Range("H126").Value = 2099.0312
$ws.Range("I126").Value = 2021.1538
$ws.Range("K126").Value = 6063.4614
$ws.Range("M126").Value = -3593.4614
$ws.Range("H136").Value = 6614.231
$ws.Range("I136").Value = 3075.8
$ws.Range("K136").Value = 9227.400000000001
$ws.Range("M136").Value = -6677.400000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 601
$ws.Range("I5").Value = 599.5
$ws.Range("K5").Value = 1798.5
$ws.Range("M5").Value = -1686.5
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()
$ws.Range("H113").Value = 2778.5454
$ws.Range("J113").Value = 2793.3333
$ws.Range("L113").Value = 8379.999899999999
$ws.Range("N113").Value = -12719.9999
$ws.Range("H131").Value = 4256795
$ws.Range("I131").Value = 7143429
$ws.Range("J131").Value = 3032162.5
$ws.Range("K131").Value = 21430287
$ws.Range("L131").Value = 9096487.5
$ws.Range("M131").Value = -21425247
$ws.Range("N131").Value = -9106567.5
$ws.Range("H135").Value = 601
$ws.Range("I135").Value = 599.5
$ws.Range("K135").Value = 5395.5
$ws.Range("M135").Value = -2860.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H99").Value = 2779.2222
$ws.Range("I99").Value = 1876.75
$ws.Range("J99").Value = 9999
$ws.Range("K99").Value = 1876.75
$ws.Range("L99").Value = 9999
$ws.Range("M99").Value = 369.25
$ws.Range("N99").Value = -14491
$ws.Range("H107").Value = 1342.75
$ws.Range("I107").Value = 1724.75
$ws.Range("K107").Value = 1724.75
$ws.Range("M107").Value = 195.25
$ws.Range("H122").Value = 3484.889
$ws.Range("I122").Value = 3487.6667
$ws.Range("J122").Value = 3479.3333
$ws.Range("K122").Value = 10463.0001
$ws.Range("L122").Value = 10437.9999
$ws.Range("M122").Value = -8013.000100000001
$ws.Range("N122").Value = -15337.9999
$ws.Range("H139").Value = 96143.23
$ws.Range("J139").Value = 96143.23
$ws.Range("L139").Value = 96143.23
$ws.Range("N139").Value = -106423.23

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 3218.3076
$ws.Range("J100").Value = 3898.5
$ws.Range("L100").Value = 3898.5
$ws.Range("N100").Value = -4980.5
$ws.Range("H102").Value = 49894
$ws.Range("J102").Value = 49894
$ws.Range("L102").Value = 49894
$ws.Range("N102").Value = -56384
$ws.Range("H110").Value = 36924.5
$ws.Range("J110").Value = 36924.5
$ws.Range("L110").Value = 36924.5
$ws.Range("N110").Value = -45104.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 5334529
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").ClearContents()
$ws.Range("H98").Value = 50000
$ws.Range("J98").Value = 50000
$ws.Range("L98").Value = 50000
$ws.Range("N98").Value = -55990
$ws.Range("H111").Value = 50643
$ws.Range("J111").Value = 50643
$ws.Range("L111").Value = 50643
$ws.Range("N111").Value = -58823
$ws.Range("H123").Value = 56799.332
$ws.Range("J123").Value = 56799.332
$ws.Range("L123").Value = 56799.332
$ws.Range("N123").Value = -66599.33199999999
$ws.Range("H132").Value = 1474.0625
$ws.Range("I132").Value = 1455.625
$ws.Range("J132").Value = 1529.375
$ws.Range("K132").Value = 4366.875
$ws.Range("L132").Value = 4588.125
$ws.Range("M132").Value = -1836.875
$ws.Range("N132").Value = -9648.125

Write-Host "Refreshed 219 market-data cells across $($wb.Worksheets.Count) sheets."
